$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.04774544758327767, 0.4624696470792148, 0.8033151732794445, 0.8962785132309289, 0.9287915074979147)
    3  = @(-0.04718447073015461, 0.5635575896698916, 0.8690252791348545, 0.9322152536484557, 0.9690367108258421)
    4  = @(-0.03314431376102876, 0.5287269987382343, 0.6606252920084058, 0.8127885899841396, 0.8482238639804189)
    5  = @(-0.1426010309797199, 0.5484211637281492, 0.6964982596166069, 0.8345647126596037, 0.862426533762813)
    6  = @(-0.2189447014214337, 0.5233672663903658, 0.6294409302107248, 0.7933731342884789, 0.8038132369664293)
    7  = @(0.05889770248139205, 0.3795808076999097, 0.2585567114336377, 0.5084847209441379, 0.5356993033274895)
    8  = @(-0.01074580271062597, 0.4136063680922605, 0.2412788381939795, 0.4912014232409954, 0.5379554248278389)
    9  = @(0.1748454189159413, 0.4191804083912449, 0.2366906289677956, 0.4865086114014793, 0.556039173689999)
    10 = @(0.7444949522781628, 0.7444949522781628, 0.5542727339676639, 0.7444949522781628)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $col = 2
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col++
    }
}
